# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column D (cultivo-descripcion), H (grupo-cultivo-descripcion) and J (secanoregadio)
# move from "dimension" metadata to "measure" metadata (and lose their mapping file).
# Column K (municipio-nombre) moves from "measure" metadata to a "dimension"/refArea
# metadata, gaining a URI-Municipio marker (like provincia-nombre / comarca-nombre).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: cultivo-descripcion ---
$ws.Range("D2").Value = "iaest-measure:cultivo-descripcion"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("D5").Clear()

# --- Column H: grupo-cultivo-descripcion ---
$ws.Range("H2").Value = "iaest-measure:grupo-cultivo-descripcion"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("H5").Clear()

# --- Column J: secanoregadio ---
$ws.Range("J2").Value = "iaest-measure:secanoregadio"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("J5").Clear()

# --- Column K: municipio-nombre ---
$ws.Range("K2").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "dim"
$ws.Range("K4").Value = "URI-Municipio"
